$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '57.117.19'
$ws.Range('E2').Value = '  -1.80%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.067.89'
$ws.Range('E3').Value = '  -1.86%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '521.36'
$ws.Range('E5').Value = '  -1.47%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '135.40'
$ws.Range('E6').Value = '  -4.91%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.066.68'
$ws.Range('E8').Value = '  -1.85%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.470'
$ws.Range('E9').Value = '  +5.58%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.31'
$ws.Range('E10').Value = '  +2.14%  '
$ws.Range('E11').Value = '  -2.82%  '
$ws.Range('E12').Value = '  +2.31%  '
$ws.Range('E13').Value = '  +1.16%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.591.40'
$ws.Range('E14').Value = '  -1.93%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '25.13'
$ws.Range('E15').Value = '  -1.96%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000160'
$ws.Range('E16').Value = '  -3.05%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '57.194.78'
$ws.Range('E17').Value = '  -1.79%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.068.94'
$ws.Range('E18').Value = '  -2.11%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.86'
$ws.Range('E19').Value = '  -4.05%  '
$ws.Range('E20').Value = '  -2.98%  '
$ws.Range('E21').Value = '  -2.02%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '348.54'
$ws.Range('E22').Value = '  +1.88%  '
$ws.Range('E23').Value = '  -0.07%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '69.10'
$ws.Range('E24').Value = '  +2.09%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.497'
$ws.Range('E26').Value = '  -2.60%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.997'
$ws.Range('E27').Value = '  -0.34%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0₃0859'
$ws.Range('E28').Value = '  -7.54%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.15'
$ws.Range('E30').Value = '  -2.89%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.85'
$ws.Range('E31').Value = '  -1.33%  '
$ws.Range('B32').Value = 'RenderToken'
$ws.Range('C32').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.79'
$ws.Range('E32').Value = '  -9.46%  '
$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '20.91'
$ws.Range('E33').Value = '  -0.89%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '158.87'
$ws.Range('E34').Value = '  +0.25%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.81'
$ws.Range('E35').Value = '  +1.20%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.12'
$ws.Range('E36').Value = '  -4.91%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.98'
$ws.Range('E37').Value = '  -3.46%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '25.21'
$ws.Range('E38').Value = '  -4.16%  '
$ws.Range('E39').Value = '  -1.69%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0653'
$ws.Range('E40').Value = '  -1.86%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.02'
$ws.Range('E41').Value = '  +0.97%  '
$ws.Range('E42').Value = '  -6.33%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.692'
$ws.Range('E43').Value = '  -0.77%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.406.17'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '36.45'
$ws.Range('E45').Value = '  -0.50%  '
$ws.Range('E46').Value = '  -0.06%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.109.45'
$ws.Range('E47').Value = '  -1.78%  '
$ws.Range('E48').Value = '  -0.78%  '
$ws.Range('E49').Value = '  -2.42%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.931'
$ws.Range('E50').Value = '  -7.06%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '19.44'
$ws.Range('E51').Value = '  -5.80%  '
